$d = $word.ActiveDocument

# 1) "IP principale" paragraph: mask /16 -> /24 (keep existing run split intact)
$para1Xml = '<w:p><w:r><w:t>IP principale : 192.</w:t></w:r><w:r><w:t>1</w:t></w:r><w:r><w:t>68.4</w:t></w:r><w:r><w:t>4.0</w:t></w:r><w:r><w:t>.</w:t></w:r><w:r><w:t>255.</w:t></w:r><w:r><w:t>255.255.0/</w:t></w:r><w:r><w:t>24</w:t></w:r></w:p>'
$d.Paragraphs.Item(1).Range.InsertXML($para1Xml)

# 2) "SR 1 : 240 PC" / "SR 2 : 144 PC" / "SR 3 : x PC" paragraphs get re-authored
#    (en-US language tagging + grammar-checker proofErr bracketing around "N :")
$sr1Xml = '<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>SR</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>1</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t> :</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> 240 PC</w:t></w:r></w:p>'
$sr2Xml = '<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>SR</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>2</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t> :</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> 144 PC</w:t></w:r></w:p>'
$sr3Xml = '<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>SR</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>3</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t> :</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>x PC</w:t></w:r></w:p>'

$d.Paragraphs.Item(3).Range.InsertXML($sr1Xml)
$d.Paragraphs.Item(4).Range.InsertXML($sr2Xml)
$d.Paragraphs.Item(5).Range.InsertXML($sr3Xml)

# 3) Move the "_GoBack" bookmark from the table cell ("...192.168.34.131") to the
#    end of the "SR non-connecté ... + 1 serveur + 1 switch" paragraph.
#    Remove it from the table cell first so only a single "_GoBack" exists at a time.
$tableCellXml = '<w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>192.168.</w:t></w:r><w:r><w:t>34</w:t></w:r><w:r><w:t>.1</w:t></w:r><w:r><w:t>31</w:t></w:r></w:p>'
$d.Paragraphs.Item(38).Range.InsertXML($tableCellXml)

$nonConnecteXml = '<w:p><w:r><w:t>SR non-connecté : 12</w:t></w:r><w:r><w:t>8 PC</w:t></w:r><w:r><w:t xml:space="preserve"> + 1 serveur + 1 switch</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
$d.Paragraphs.Item(6).Range.InsertXML($nonConnecteXml)
